# Bootstrap metrics/SOTA models added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("dt" strategy) - updated bootstrap metric means/stds and CIs
$ws.Range("B2").Value = 0.959175717070454
$ws.Range("C2").Value = 0.899444958371878
$ws.Range("E2").Value = 0.8706204236080397
$ws.Range("F2").Value = 0.9047619047619048
$ws.Range("G2").Value = 0.9242656449553002
$ws.Range("H2").Value = 0.02151508790621114
$ws.Range("I2").Value = 0.03806435376428949
$ws.Range("K2").Value = 0.03673459897686258
$ws.Range("L2").Value = 0.03072531109379957
$ws.Range("M2").Value = 0.05647417221222185
$ws.Range("N2").Value = "[0.944, 0.975]"
$ws.Range("O2").Value = "[0.872, 0.927]"
$ws.Range("Q2").Value = "[0.844, 0.897]"
$ws.Range("R2").Value = "[0.883, 0.927]"

# Row 3 ("knn" strategy) - updated roc_auc_p
$ws.Range("T3").Value = 0.06862723695369878

# Row 4 ("linear" strategy) - updated roc_auc_p
$ws.Range("T4").Value = 0.1366846762177055

# Row 5 ("simple" strategy) - updated roc_auc_p
$ws.Range("T5").Value = 0.3674628466525119
